$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 keeps referencing the same shared string slot, but its text changes from "asd" to "Julitka"
$ws.Range("A1").Value = "Julitka"

# B1 previously shared the same string as A1/C1 ("asd"); it now becomes its own new string "Doma}ska"
$ws.Range("B1").Value = "Doma}ska"

# C1 becomes a numeric barcode value, formatted like the other barcode cells in column C (style index 1)
$ws.Range("C1").Value = 300621357439
$ws.Range("C1").NumberFormat = "0"

# Update the active selection/cursor cell to I4, matching the saved view state
$ws.Range("I4").Select()
